$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at row 204, shifting existing rows 204..292 down to 205..293
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new data record
$ws.Cells.Item(204, 1).Value = 7
$ws.Cells.Item(204, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(204, 3).Value = "Ñuble"
$ws.Cells.Item(204, 4).Value = 44875
$ws.Cells.Item(204, 5).Value = 16
$ws.Cells.Item(204, 6).Value = 100112003
$ws.Cells.Item(204, 7).Value = "Ajo"
$ws.Cells.Item(204, 8).Value = "Chino"
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 60
$ws.Cells.Item(204, 11).Value = 17000
$ws.Cells.Item(204, 12).Value = 18000
$ws.Cells.Item(204, 13).Value = 17500
$ws.Cells.Item(204, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(204, 15).Value = "China"
$ws.Cells.Item(204, 16).Value = 1750
$ws.Cells.Item(204, 17).Value = 10
$ws.Cells.Item(204, 18).Value = "Hortaliza"
